$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
